# Applies the edit described by the commit "Added many more features" to
# the Griffin's Quest review document:
#   - New title / recap heading text (2 occurrences)
#   - Rewritten "What we like" bullet list (3 bullets)
#   - Rewritten "What we don't like" bullet list (2 bullets)
#   - New meta description (italic paragraph at the end)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the (single) paragraph whose visible text equals $exactText
# (paragraph/cell marks stripped).
# ---------------------------------------------------------------------------
function Find-ParagraphsByText($exactText) {
    $results = @()
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $exactText) {
            $results += $p
        }
    }
    return $results
}

# ---------------------------------------------------------------------------
# Replace the text of a paragraph that holds it in a single, already-
# formatted run (e.g. the bold recap heading or the italic meta description).
# Uses a plain Range.Text assignment (instead of Find.Execute) so that no
# "smart quote" autocorrection is applied to the replacement text - this
# lets us control straight vs curly apostrophes exactly.
# ---------------------------------------------------------------------------
function Set-ParagraphText($oldText, $newText) {
    $paras = Find-ParagraphsByText $oldText
    foreach ($para in $paras) {
        $rng = $d.Range($para.Range.Start, $para.Range.Start + $oldText.Length)
        $rng.Text = $newText
    }
}

# ---------------------------------------------------------------------------
# Replace the text of a ListBullet paragraph that starts with an empty
# <w:r/> run followed by the text run. A direct Range.Text assignment (or a
# Find.Execute replace) causes this runtime to silently drop that leading
# empty run, so instead we rebuild the whole paragraph via InsertXML,
# re-using its original <w:pPr> and keeping the leading empty run intact.
# ---------------------------------------------------------------------------
function Set-BulletParagraphText($oldText, $newText) {
    $paras = Find-ParagraphsByText $oldText
    foreach ($para in $paras) {
        $openXml = $para.Range.WordOpenXML
        $bodyStart = $openXml.IndexOf("<w:body>")
        $pStart = $openXml.IndexOf("<w:p", $bodyStart)
        $pOpenEnd = $openXml.IndexOf(">", $pStart) + 1
        $pPrXml = ""
        $pPrStart = $openXml.IndexOf("<w:pPr>", $pOpenEnd)
        if ($pPrStart -ge 0 -and $pPrStart -lt ($pOpenEnd + 20)) {
            $pPrEnd = $openXml.IndexOf("</w:pPr>", $pPrStart) + "</w:pPr>".Length
            $pPrXml = $openXml.Substring($pPrStart, $pPrEnd - $pPrStart)
        }
        $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $pPrXml + '<w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p>'
        $para.Range.InsertXML($newParaXml)
    }
}

# Title heading (Heading1) and the bold recap line further down share the
# exact same old/new text, so a single call handles both occurrences.
Set-ParagraphText "Play Griffin's Quest Free Slot - Review of Kalamba Games' Slot Game" ("Play Griffin" + [char]0x2019 + "s Quest for Free - Review of Gameplay and Features")

# "What we like" bullet list
Set-BulletParagraphText "Unique game mechanism with a high number of pay lines" "Unique game mechanism with variable-length reels"
Set-BulletParagraphText "High RTP of 97.52%" "Large number of pay lines for more winning opportunities"
Set-BulletParagraphText "Varied gameplay with a bonus mode for big wins" "Exciting bonus features"

# "What we don't like" bullet list
Set-BulletParagraphText "Very high volatility requires a lot of patience and a generous bankroll" "Very high volatility requiring patience and a generous bankroll"
Set-BulletParagraphText "Gameplay is better suited for long gaming sessions" "Not suitable for players looking for quick gameplay"

# Meta description (italic paragraph at the very end). Its new text keeps a
# straight apostrophe in "Griffin's" (unlike the title text above).
Set-ParagraphText "Read our review of the Griffin's Quest slot game by Kalamba Games. Play now and enjoy the high-quality design, unique game mechanism, and varied gameplay with numerous chances for big wins. Free demo available." ("Discover the highlights of Griffin" + [char]39 + "s Quest online slot game and play for free. Learn about the gameplay, bonus features, and more.")

Write-Output "Edits applied."
